# Workbook / sheet handles -----------------------------------------------
$wb = $excel.ActiveWorkbook

$wsNodes  = $wb.Worksheets.Item("config_gridNodes")
$wsConn   = $wb.Worksheets.Item("config_gridConnections")
$wsAssets = $wb.Worksheets.Item("config_energyAssets")

# --------------------------------------------------------------------
# 1. Insert 4 new rows in config_energyAssets (pushes the old row 5
#    -- id a4 / PRODUCTION / Solarpanels_1MW / b3 -- down to row 9,
#    and gives rows 6-9 the same formatting as row 5 had).
# --------------------------------------------------------------------
$wsAssets.Rows("6:9").Insert()

# --------------------------------------------------------------------
# 2. Row 5 becomes a brand new STORAGE/"EHGV" asset (id a4, parent b2).
#    Typing fresh values into these cells drops their old style (4),
#    matching the un-styled look of the newly-entered data.
# --------------------------------------------------------------------
$wsAssets.Range("A5").ClearFormats()
$wsAssets.Range("A5").Value2 = 3
$wsAssets.Range("C5").ClearFormats()
$wsAssets.Range("C5").Value2 = "a4"
$wsAssets.Range("D5").ClearFormats()
$wsAssets.Range("D5").Value2 = "STORAGE"
$wsAssets.Range("E5").ClearFormats()
$wsAssets.Range("E5").Value2 = "EHGV"
$wsAssets.Range("F5").ClearFormats()
$wsAssets.Range("F5").Value2 = "b2"

# --------------------------------------------------------------------
# 3. Rows 6-8: three more new STORAGE / "EHGV" assets (a5, a6, a7),
#    all parented on b2. These rows inherited style 4 from the Insert,
#    and only A/C keep being typed without it (same pattern as row 5).
# --------------------------------------------------------------------
$newAssets = @(
    @{ Row = 6; Index = 4; Id = "a5" },
    @{ Row = 7; Index = 5; Id = "a6" },
    @{ Row = 8; Index = 6; Id = "a7" }
)

foreach ($a in $newAssets) {
    $r = $a.Row
    $wsAssets.Range("A$r").ClearFormats()
    $wsAssets.Range("A$r").Value2 = $a.Index
    $wsAssets.Range("B$r").Value2 = "energyAsset"
    $wsAssets.Range("C$r").ClearFormats()
    $wsAssets.Range("C$r").Value2 = $a.Id
    $wsAssets.Range("D$r").Value2 = "STORAGE"
    $wsAssets.Range("E$r").Value2 = "EHGV"
    $wsAssets.Range("F$r").Value2 = "b2"
}

# --------------------------------------------------------------------
# 4. Row 9 keeps the original row-5 content (id a4 -> renamed a8 now
#    that it is the 8th asset, index bumped to 7). Style 4 carried
#    through the row insert, so it is left untouched.
# --------------------------------------------------------------------
$wsAssets.Range("A9").Value2 = 7
$wsAssets.Range("B9").Value2 = "energyAsset"
$wsAssets.Range("C9").Value2 = "a8"
$wsAssets.Range("D9").Value2 = "PRODUCTION"
$wsAssets.Range("E9").Value2 = "Solarpanels_1MW"
$wsAssets.Range("F9").Value2 = "b3"

# --------------------------------------------------------------------
# 5. Selections / active sheet, mirroring what a user would leave
#    behind after this editing session. config_energyAssets ends up
#    both the active sheet and the one with focus (tabSelected).
# --------------------------------------------------------------------
$wsNodes.Select() | Out-Null
$wsNodes.Range("I15").Select() | Out-Null

$wsConn.Select() | Out-Null
$wsConn.Range("C42").Select() | Out-Null

$wsAssets.Select() | Out-Null
$wsAssets.Range("D14").Select() | Out-Null
